$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together as a block for each data row (2..16)
$cols = @("D","J","K","L","M","N","O","P","Q")

# Snapshot the current (pre-edit) values for the columns above, per row,
# before any writes happen, so the permutation below is safe regardless
# of write order.
$snapshot = @{}
for ($r = 2; $r -le 16; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Target row -> source row (data originating from source row now belongs
# to target row), derived from the diff.
$mapping = @{
    2  = 6
    3  = 15
    4  = 10
    5  = 13
    6  = 8
    7  = 7
    8  = 11
    9  = 2
    10 = 5
    11 = 16
    12 = 14
    13 = 3
    14 = 4
    15 = 12
    16 = 9
}

foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    $srcVals = $snapshot[$sourceRow]
    foreach ($c in $cols) {
        $ws.Range("$c$targetRow").Value = $srcVals[$c]
    }
}
